$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Recommandations): header row 1 is unchanged ---
# Data rows 2-42, columns A-G
$sheet1Data = @(
    @('BRVM - SERVICES PUBLICS', 0, 8, 3382.68, 115.04, '🟡 Observer', '➖ Neutre'),
    @('AIR LIQUIDE CI', 0, 4, 2775, 690, '🟡 Observer', '➖ Neutre'),
    @('NEI-CEDA CI', 0, 3, 2545, 890, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AUTRES SECTEURS', 0, 4, 2382.43, 589.61, '🟡 Observer', '➖ Neutre'),
    @('BRVM - DISTRIBUTION', 0, 4, 2257.22, 525.41, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TRANSPORT', 0, 4, 1494.34, 396.53, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AGRICULTURE', 0, 4, 1313.55, 315.23, '🟡 Observer', '➖ Neutre'),
    @('SETAO CI', 0, 1, 995, 995, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 4, 840.33, 192.5, '🟡 Observer', '➖ Neutre'),
    @('BRVM - FINANCES', 0, 4, 568.64, 141.46, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRESTIGE', 0, 4, 565.01, 139.5, '🟡 Observer', '➖ Neutre'),
    @('BRVM - SERVICES FINANCIERS', 0, 4, 558.85, 139.02, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIELS', 0, 4, 488.3, 121.84, '🟡 Observer', '➖ Neutre'),
    @('BRVM - ENERGIE', 0, 4, 424.38, 104.38, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TELECOMMUNICATIONS', 0, 4, 380.33, 94.48, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIE                  (**)', 0, 1, 229.16, 229.16, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRINCIPAL                   (**)', 0, 1, 206.91, 206.91, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DE BASE         (**)', 0, 1, 194.06, 194.06, '🟡 Observer', '➖ Neutre'),
    @('NEI-CEDA CI (NEIC)', 4, 0, 29.12, 7.37, '🟢 Achat', '✅ Renforcer'),
    @('SOLIBRA CI (SLBC)', 2, 0, 14.98, 7.49, '🟡 Observer', '➖ Neutre'),
    @('SICABLE CI (CABC)', 2, 0, 14.55, 7.3, '🟡 Observer', '➖ Neutre'),
    @('AFRICA GLOBAL LOGISTICS CI (SDSC)', 3, 0, 9.94, 6.23, '🟢 Achat', '✅ Renforcer'),
    @('SICOR CI (SICC)', 1, 0, 6.06, 6.06, '🟡 Observer', '➖ Neutre'),
    @('ECOBANK TRANS. INCORP. TG (ETIT)', 1, 0, 4.55, 4.55, '🟡 Observer', '➖ Neutre'),
    @('ONATEL BF (ONTBF)', 1, 0, 2.73, 2.73, '🟡 Observer', '➖ Neutre'),
    @('SOCIETE IVOIRIENNE DE BANQUE  (SIBC)', 1, 0, 2.59, 2.59, '🟡 Observer', '➖ Neutre'),
    @('SITAB CI (STBC)', 1, 0, 2.56, 2.56, '🟡 Observer', '➖ Neutre'),
    @('SAPH CI (SPHC)', 1, 0, 2.53, 2.53, '🟡 Observer', '➖ Neutre'),
    @('ECOBANK COTE D''''IVOIRE (ECOC)', 1, 0, 1.97, 1.97, '🟡 Observer', '➖ Neutre'),
    @('SERVAIR ABIDJAN CI (ABJC)', 1, 0, 1.92, 1.92, '🟡 Observer', '➖ Neutre'),
    @('TOTAL', 0, 4, 0, 0, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA ML (BOAM)', 1, 1, -2.82, 4.65, '🟡 Observer', '👀 À surveiller'),
    @('SUCRIVOIRE (SCRC)', 0, 1, -7.36, -7.36, '🟡 Observer', '➖ Neutre'),
    @('SETAO CI (STAC)', 0, 1, -7.44, -7.44, '🟡 Observer', '➖ Neutre'),
    @('CFAO MOTORS CI (CFAC)', 0, 1, -7.48, -7.48, '🟡 Observer', '➖ Neutre'),
    @('TRACTAFRIC MOTORS CI (PRSC)', 0, 2, -14.92, -7.42, '🟡 Observer', '➖ Neutre'),
    @('PALM CI (PALC)', 0, 2, -14.93, -7.44, '🟡 Observer', '➖ Neutre'),
    @('FILTISAC CI (FTSC)', 0, 2, -14.98, -7.48, '🟡 Observer', '➖ Neutre'),
    @('SAFCA CI (SAFC)', 0, 3, -22.32, -7.36, '🔴 Vente', '⚠️ Risque de décrochage'),
    @('BERNABE CI (BNBC)', 0, 3, -22.41, -7.49, '🔴 Vente', '⚠️ Risque de décrochage'),
    @('NESTLE CI (NTLC)', 0, 4, -29.92, -7.49, '🔴 Vente', '⚠️ Risque de décrochage')
)

$r = 2
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# --- Sheet2 (Top_YTD): header row 1 is unchanged ---
# Data rows 2-11, columns A-B
$sheet2Data = @(
    @('BRVM - SERVICES PUBLICS', 10098391.15),
    @('AIR LIQUIDE CI', 396827.6),
    @('BRVM - AUTRES SECTEURS', 233990.61),
    @('BRVM - DISTRIBUTION', 193846.51),
    @('NEI-CEDA CI', 84738.05),
    @('BRVM - TRANSPORT', 50115.01),
    @('BRVM - AGRICULTURE', 33523.84),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 9107.67),
    @('BRVM - FINANCES', 3338.74),
    @('BRVM-PRESTIGE', 3287.27)
)

$r = 2
foreach ($row in $sheet2Data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r++
}
